$wb = $excel.ActiveWorkbook

# --- process_params sheet ---
$ws1 = $wb.Worksheets.Item("process_params")

# New row: C7 = C4 - B4 (difference used for the updated Scaling "Threshold")
$ws1.Range("C7").Formula = "=C4-B4"

# Leave the cursor where the author left it after entering the formula
$ws1.Range("C8").Select()

# --- Scaling sheet ---
$ws2 = $wb.Worksheets.Item("Scaling")

# Weighting update: clear out the old Min/Max/Inversion numbers (B:D) for
# rows 2-4, keeping their existing cell formatting/style intact.
$ws2.Range("B2:D4").ClearContents()

# Scaling update: new Threshold value for PLA_virgin (F4), matching the
# newly computed C4-B4 difference on process_params.
$ws2.Range("F4").Value = 25.332000000000001

# Selection ends up on F5 after editing F4
$ws2.Range("F5").Select()
